# This document is a "two-digit divided by one-digit" practice sheet:
# a dated heading paragraph followed by a 5-column table whose data rows
# (1, 5, 9, 13, 17 -- the others are blank spacer rows) hold the division
# problems. The edit bumps the heading date by one day and swaps every
# problem/answer string for a new one.
$d = $word.ActiveDocument

function Replace-InRange([int]$rangeStart, [int]$rangeEnd, [string]$old, [string]$new) {
    # Scope the replacement to an explicit Document.Range built from fixed
    # character offsets and use wdReplaceOne (1): Find.Execute on this host
    # can otherwise walk past the end of a loosely-held Range (e.g. a
    # Cell.Range reference) and edit a later occurrence of identical text
    # sitting in a different table cell, which matters here because several
    # cells share the same original text but need different replacements.
    $scoped = $d.Range($rangeStart, $rangeEnd)
    $ok = $scoped.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 1)
    if (-not $ok) {
        throw ("Could not find expected text " + $old + " in range " + $rangeStart + "-" + $rangeEnd)
    }
}

# Heading date line.
Replace-InRange $d.Content.Start $d.Content.End "2025-12-23 Tuesday" "2025-12-24 Wednesday"

# Table data cells, addressed by (row, column) so the handful of duplicate
# original strings (e.g. "61÷3=20, 1" and "77÷8=9, 5" each occur twice)
# resolve to the correct, cell-specific replacement.
$t = $d.Tables.Item(1)
$cell = $t.Cell(1, 1)
Replace-InRange $cell.Range.Start $cell.Range.End "43÷8=5, 3" "17÷9=1, 8"
$cell = $t.Cell(1, 2)
Replace-InRange $cell.Range.Start $cell.Range.End "84÷6=14, 0" "96÷5=19, 1"
$cell = $t.Cell(1, 3)
Replace-InRange $cell.Range.Start $cell.Range.End "62÷8=7, 6" "31÷4=7, 3"
$cell = $t.Cell(1, 4)
Replace-InRange $cell.Range.Start $cell.Range.End "18÷6=3, 0" "56÷9=6, 2"
$cell = $t.Cell(1, 5)
Replace-InRange $cell.Range.Start $cell.Range.End "30÷5=6, 0" "63÷9=7, 0"
$cell = $t.Cell(5, 1)
Replace-InRange $cell.Range.Start $cell.Range.End "11÷3=3, 2" "26÷2=13, 0"
$cell = $t.Cell(5, 2)
Replace-InRange $cell.Range.Start $cell.Range.End "43÷2=21, 1" "38÷3=12, 2"
$cell = $t.Cell(5, 3)
Replace-InRange $cell.Range.Start $cell.Range.End "61÷3=20, 1" "38÷9=4, 2"
$cell = $t.Cell(5, 4)
Replace-InRange $cell.Range.Start $cell.Range.End "44÷6=7, 2" "24÷6=4, 0"
$cell = $t.Cell(5, 5)
Replace-InRange $cell.Range.Start $cell.Range.End "10÷6=1, 4" "45÷6=7, 3"
$cell = $t.Cell(9, 1)
Replace-InRange $cell.Range.Start $cell.Range.End "92÷3=30, 2" "37÷7=5, 2"
$cell = $t.Cell(9, 2)
Replace-InRange $cell.Range.Start $cell.Range.End "61÷3=20, 1" "19÷8=2, 3"
$cell = $t.Cell(9, 3)
Replace-InRange $cell.Range.Start $cell.Range.End "27÷4=6, 3" "60÷6=10, 0"
$cell = $t.Cell(9, 4)
Replace-InRange $cell.Range.Start $cell.Range.End "99÷5=19, 4" "92÷6=15, 2"
$cell = $t.Cell(9, 5)
Replace-InRange $cell.Range.Start $cell.Range.End "45÷4=11, 1" "33÷2=16, 1"
$cell = $t.Cell(13, 1)
Replace-InRange $cell.Range.Start $cell.Range.End "85÷2=42, 1" "54÷2=27, 0"
$cell = $t.Cell(13, 2)
Replace-InRange $cell.Range.Start $cell.Range.End "46÷8=5, 6" "77÷8=9, 5"
$cell = $t.Cell(13, 3)
Replace-InRange $cell.Range.Start $cell.Range.End "60÷4=15, 0" "62÷5=12, 2"
$cell = $t.Cell(13, 4)
Replace-InRange $cell.Range.Start $cell.Range.End "46÷3=15, 1" "76÷5=15, 1"
$cell = $t.Cell(13, 5)
Replace-InRange $cell.Range.Start $cell.Range.End "62÷7=8, 6" "53÷5=10, 3"
$cell = $t.Cell(17, 1)
Replace-InRange $cell.Range.Start $cell.Range.End "77÷8=9, 5" "96÷8=12, 0"
$cell = $t.Cell(17, 2)
Replace-InRange $cell.Range.Start $cell.Range.End "40÷2=20, 0" "65÷2=32, 1"
$cell = $t.Cell(17, 3)
Replace-InRange $cell.Range.Start $cell.Range.End "54÷6=9, 0" "49÷6=8, 1"
$cell = $t.Cell(17, 4)
Replace-InRange $cell.Range.Start $cell.Range.End "91÷2=45, 1" "20÷3=6, 2"
$cell = $t.Cell(17, 5)
Replace-InRange $cell.Range.Start $cell.Range.End "81÷4=20, 1" "87÷2=43, 1"

Write-Host "Applied 26 text replacements (1 heading + 25 table cells)."
